# Fixed state social benefits add factor
# Updates recomputed projection values for rows 2-9 (state_health_outlays,
# state_social_benefits, state_non_corporate_taxes, state_corporate_taxes,
# federal_social_benefits, federal_subsidies, consumption_grants) that
# changed as a result of correcting the state social benefits add factor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2" = 171.237
    "G2" = 175.56274099593
    "H2" = 179.888481991861
    "I2" = 180.498883835544
    "J2" = 181.109285679227
    "K2" = 181.71968752291
    "L2" = 182.330089366593
    "M2" = 180.35508255879
    "N2" = 178.380075750986
    "O2" = 176.405068943183
    "P2" = 174.430062135379
    "Q2" = 175.461672281072
    "R2" = 176.493282426765
    "S2" = 177.524892572458
    "T2" = 178.556502718151
    "U2" = 181.209999767016
    "V2" = 183.86349681588
    "W2" = 186.516993864745
    "X2" = 189.170490913609
    "Y2" = 191.812994523314
    "Z2" = 194.455498133018
    "AA2" = 197.098001742723
    "AB2" = 199.740505352428
    "AC2" = 202.769995751918
    "AD2" = 205.799486151409
    "AE2" = 208.828976550899
    "AF2" = 211.85846695039
    "AG2" = 214.912515274456
    "AH2" = 217.966563598521
    "AI2" = 221.020611922587
    "AJ2" = 224.074660246653
    "AK2" = 227.230397882949
    "AL2" = 230.386135519244
    "AM2" = 233.54187315554
    "AN2" = 236.697610791836
    "AO2" = 240.027204349042
    "AP2" = 243.356797906247
    "AQ2" = 246.686391463453
    "AR2" = 250.015985020658
    "AS2" = 253.475726712435
    "AT2" = 256.935468404212
    "AU2" = 260.39521009599
    "AV2" = 263.854951787767
    "B3" = 143
    "C3" = 143.4
    "D3" = 144.8
    "E3" = 145.2
    "F3" = 145.9
    "G3" = 148.620995113912
    "H3" = 147.797019737607
    "I3" = 159.242656111224
    "J3" = 170.909193879409
    "K3" = 182.800138368449
    "L3" = 193.820513614133
    "M3" = 214.404660294731
    "N3" = 234.020541720861
    "O3" = 253.778864967589
    "P3" = 273.526701757442
    "Q3" = 281.501353023059
    "R3" = 289.895246875818
    "S3" = 298.670270353253
    "T3" = 306.911712356143
    "U3" = 308.004386827056
    "V3" = 309.058948335501
    "W3" = 310.075396881479
    "X3" = 311.091845427458
    "Y3" = 312.152874549787
    "Z3" = 313.213903672117
    "AA3" = 314.313045756914
    "AB3" = 315.412187841712
    "AC3" = 314.980134697246
    "AD3" = 314.54808155278
    "AE3" = 314.154141370781
    "AF3" = 313.79831415125
    "AG3" = 313.342899676483
    "AH3" = 312.887485201714
    "AI3" = 312.470183689414
    "AJ3" = 312.052882177113
    "AK3" = 311.261323296029
    "AL3" = 310.469764414943
    "AM3" = 309.678205533858
    "AN3" = 308.924759615241
    "AO3" = 307.542519377805
    "AP3" = 306.198392102836
    "AQ3" = 304.930490752802
    "AR3" = 303.738815327704
    "AS3" = 302.095589326862
    "AT3" = 300.528589250955
    "AU3" = 299.07592806245
    "AV3" = 297.661379836413
    "F4" = 1965.1
    "G4" = 1985.24481783359
    "H4" = 2013.98715613545
    "I4" = 2040.75997876518
    "J4" = 2066.09388583537
    "K4" = 2087.01212078221
    "L4" = 2109.46819673332
    "M4" = 2131.62749635026
    "N4" = 2152.02412440678
    "O4" = 2174.84892247001
    "P4" = 2198.85183264763
    "Q4" = 2223.5921870495
    "R4" = 2246.93859200308
    "S4" = 2271.69693284945
    "T4" = 2296.64413136301
    "U4" = 2322.6075639905
    "V4" = 2347.99543039417
    "W4" = 2373.51819513155
    "X4" = 2398.48338008961
    "Y4" = 2423.32265993621
    "Z4" = 2448.36878389449
    "AA4" = 2472.89330096243
    "AB4" = 2497.15701458521
    "AC4" = 2520.98006031788
    "AD4" = 2544.47035682741
    "AE4" = 2567.96964655918
    "AF4" = 2591.8016855141
    "AG4" = 2616.1463381371
    "AH4" = 2640.90467898348
    "AI4" = 2665.99576905299
    "AJ4" = 2691.04189301127
    "AK4" = 2716.18694241427
    "AL4" = 2741.4129308175
    "AM4" = 2766.59395310949
    "AN4" = 2791.67604995676
    "AO4" = 2816.99197058245
    "AP4" = 2842.37983698612
    "AQ4" = 2868.0734729462
    "AR4" = 2894.26173612987
    "AS4" = 2920.90865364814
    "AT4" = 2948.07717805674
    "AU4" = 2975.34462791007
    "AV4" = 3002.76496254159
    "F5" = 96.6
    "G5" = 97.5902750001144
    "H5" = 99.0031852234919
    "I5" = 100.319278382126
    "J5" = 101.564637612181
    "K5" = 102.592932098907
    "L5" = 103.696823471802
    "M5" = 104.786125971928
    "N5" = 105.788779409544
    "O5" = 106.910796351638
    "P5" = 108.090726697756
    "Q5" = 109.306908182272
    "R5" = 110.454566173476
    "S5" = 111.671631832098
    "T5" = 112.897981318847
    "U5" = 114.174286642655
    "V5" = 115.422298395032
    "W5" = 116.676941453212
    "X5" = 117.904175114069
    "Y5" = 119.125219556174
    "Z5" = 120.356432000513
    "AA5" = 121.562003395741
    "AB5" = 122.754754266415
    "AC5" = 123.925842871461
    "AD5" = 125.080574255523
    "AE5" = 126.235747726638
    "AF5" = 127.407278418738
    "AG5" = 128.604008072894
    "AH5" = 129.821073731517
    "AI5" = 131.054496611124
    "AJ5" = 132.285709055462
    "AK5" = 133.521784457391
    "AL5" = 134.761838642802
    "AM5" = 135.999682392945
    "AN5" = 137.232663185498
    "AO5" = 138.477138241446
    "AP5" = 139.725149993822
    "AQ5" = 140.988192706021
    "AR5" = 142.27555020617
    "AS5" = 143.585454146054
    "AT5" = 144.920999135047
    "AU5" = 146.261407081631
    "AV5" = 147.609330508126
    "F7" = 1805.1
    "G7" = 1626.02571860579
    "H7" = 1498.53485474677
    "I7" = 1511.08765153362
    "J7" = 1512.38995417448
    "K7" = 1506.67098308902
    "L7" = 1490.44612277592
    "M7" = 1509.63823634565
    "N7" = 1519.632602017
    "O7" = 1523.13445574117
    "P7" = 1519.17194815658
    "Q7" = 1526.85403734435
    "R7" = 1535.48425853468
    "S7" = 1543.71889409799
    "T7" = 1552.4850878848
    "U7" = 1588.56945568078
    "V7" = 1619.05458955385
    "W7" = 1642.5093049925
    "X7" = 1659.61223551309
    "Y7" = 1667.58404447224
    "Z7" = 1678.25356981618
    "AA7" = 1690.72593660438
    "AB7" = 1705.40653110591
    "AC7" = 1718.34111239517
    "AD7" = 1730.9634667895
    "AE7" = 1742.36951186696
    "AF7" = 1753.01705483799
    "AG7" = 1784.68057949284
    "AH7" = 1812.35462710207
    "AI7" = 1835.2404610672
    "AJ7" = 1853.72632303349
    "AK7" = 1827.9445084856
    "AL7" = 1814.12430186748
    "AM7" = 1812.00479390401
    "AN7" = 1821.67292333428
    "AO7" = 1878.9537705036
    "AP7" = 1923.81239873781
    "AQ7" = 1957.35253104793
    "AR7" = 1978.92921847087
    "AS7" = 1939.9220024118
    "AT7" = 1901.58652711119
    "AU7" = 1865.37884105868
    "AV7" = 1830.56104730358
    "F8" = 402.7
    "G8" = 404.587256188855
    "H8" = 406.517873843756
    "I8" = 408.48772330127
    "J8" = 410.480285907665
    "K8" = 412.460459523761
    "L8" = 414.40966066411
    "M8" = 416.34647281416
    "N8" = 418.279155300777
    "O8" = 420.218032282544
    "P8" = 422.175492749759
    "Q8" = 424.151536702422
    "R8" = 426.139969645384
    "S8" = 428.124272924913
    "T8" = 430.114770699592
    "U8" = 432.109398137703
    "V8" = 434.101960744098
    "W8" = 436.069745369896
    "X8" = 438.023076173678
    "Y8" = 439.961953155445
    "Z8" = 441.888441146912
    "AA8" = 443.796345652932
    "AB8" = 445.702185327235
    "AC8" = 447.610089833255
    "AD8" = 449.520059170991
    "AE8" = 451.432093340444
    "AF8" = 453.356516500195
    "AG8" = 455.285069323379
    "AH8" = 457.22188147343
    "AI8" = 459.183471604078
    "AJ8" = 461.16571005189
    "AK8" = 463.170661648584
    "AL8" = 465.194196730726
    "AM8" = 467.232185634883
    "AN8" = 469.268109707324
    "AO8" = 471.299904116332
    "AP8" = 473.333763357056
    "AQ8" = 475.367622597781
    "AR8" = 477.401481838505
    "AS8" = 479.441535574379
    "AT8" = 481.483654141969
    "AU8" = 483.529902372993
    "AV8" = 485.5740857723
    "F9" = 262.437
    "G9" = 257.265991927795
    "H9" = 256.215109642153
    "I9" = 258.570678761324
    "J9" = 261.162791929795
    "K9" = 263.648127833867
    "L9" = 266.14790345037
    "M9" = 268.705699022449
    "N9" = 271.263979862962
    "O9" = 273.87029090685
    "P9" = 276.524049592054
    "Q9" = 279.241851093381
    "R9" = 282.013414803142
    "S9" = 284.805174670076
    "T9" = 287.631080748397
    "U9" = 290.487699974147
    "V9" = 293.356984616877
    "W9" = 296.242478178848
    "X9" = 299.141687420886
    "Y9" = 302.064126557864
    "Z9" = 305.006516878386
    "AA9" = 307.96222864379
    "AB9" = 310.947508525384
    "AC9" = 313.960576919015
    "AD9" = 317.004049216013
    "AE9" = 320.073239044574
    "AF9" = 323.178072165758
    "AG9" = 326.312958747243
    "AH9" = 329.463709405131
    "AI9" = 332.649728371015
    "AJ9" = 335.868341966415
    "AK9" = 339.121324982229
    "AL9" = 342.398384134518
    "AM9" = 345.704227507039
    "AN9" = 349.026833190491
    "AO9" = 352.367229319003
    "AP9" = 355.72749945931
    "AQ9" = 359.106148830796
    "AR9" = 362.503210437325
    "AS9" = 365.926019735059
    "AT9" = 369.368953535823
    "AU9" = 372.838895027704
    "AV9" = 376.325951842415
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
